$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ridership")

# Update the Riders (column C) and Average (column D) figures for the
# ridership run on 20161026.
$ws.Range("C2").Value = 183
$ws.Range("D2").Value = 104.09

$ws.Range("C3").Value = 240
$ws.Range("D3").Value = 109.64

$ws.Range("C4").Value = 227
$ws.Range("D4").Value = 108.23

$ws.Range("C5").Value = 210
$ws.Range("D5").Value = 104.9

$ws.Range("D6").Value = 49.59

$ws.Range("C7").Value = 73
$ws.Range("D7").Value = 38.13
